$d = $word.ActiveDocument

# Locate the run "A user may not be bothered with escaping:" and note its
# start position so the text can be split into two runs, matching the
# author's edit (a reminder about Python 3 support was appended).
$rng = $d.Content
$rng.Find.Execute("A user may not be bothered with escaping:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Split point: keep "A user ma" in the first run, replace the remainder.
$splitAt = $rng.Start + 9
$tail = $d.Range($splitAt, $rng.End)
$tail.Text = "y not be bothered with escaping. Note it works only for Python 3:"

# Toggling then restoring a character property forces Word to keep the
# replaced text as its own run instead of re-merging it with the
# preceding run, producing the same two-run structure as the edit.
$tail.Bold = 1
$tail.Bold = 0
